# PLS 1month with new parameters
# Append a new data row (row 3) for the PLS model to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These three columns hold numeric-looking values that must be stored as
# TEXT (matching columns F2/G2/H2 etc. in the existing sheet), not as
# numbers. Use a leading apostrophe to force text entry, then reset the
# cell style back to Normal so no stray number-format style lingers.
$ws.Range("C3").Value = "'0.4315"
$ws.Range("C3").Style = "Normal"

$ws.Range("F3").Value = "'0.3538"
$ws.Range("F3").Style = "Normal"

$ws.Range("H3").Value = "'0.2779"
$ws.Range("H3").Style = "Normal"

# --- Text cells -------------------------------------------------------
$ws.Range("A3").Value = "PLS"
$ws.Range("B3").Value = "PLS on data with 1h sampling over 1 month"

# --- Numeric cells ------------------------------------------------------
$ws.Range("D3").Value = -0.4345
$ws.Range("E3").Value = 0
$ws.Range("G3").Value = -0.1648
$ws.Range("I3").Value = -0.13689999999999999
$ws.Range("J3").Value = -0.1842
$ws.Range("K3").Value = 0.58599999999999997

# --- Selection moved to B8 as in the edited workbook ---------------------
$ws.Range("B8").Select()
